$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) updates keep their original text representation
# (values like "0.470", "7.70", "61.116.27" must not be auto-converted to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.116.27'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.397.51'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.03'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.82'
$ws.Range("E6").Value = '  +1.09%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.397.13'
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.70'
$ws.Range("E10").Value = '  +2.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.122'
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.380'
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.981.88'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("E14").Value = '  -0.97%  '
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  -2.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.400.33'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.215.02'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.95'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("E20").Value = '  -1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.31'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '375.45'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.531.32'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.552'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.16'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.176'
$ws.Range("E28").Value = '  +9.30%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.61'
$ws.Range("E29").Value = '  -6.78%  '
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.41'
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.13'
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.44'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("E36").Value = '  +1.80%  '
$ws.Range("E37").Value = '  -3.06%  '
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '166.15'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0772'
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.85'
$ws.Range("E41").Value = '  +5.79%  '
$ws.Range("E42").Value = '  +0.72%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.516.91'
$ws.Range("E48").Value = '  +6.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.85'
$ws.Range("E49").Value = '  +5.36%  '
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("E51").Value = '  -0.29%  '
